$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.464.40"
$ws.Range("E2").Value = "  +1.61%  "
$ws.Range("D3").Value = "1.676.28"
$ws.Range("E3").Value = "  +2.40%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'217.28"
$ws.Range("E5").Value = "  +1.66%  "
$ws.Range("D6").Value = "'0.5320"
$ws.Range("E6").Value = "  +1.48%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.2693"
$ws.Range("E8").Value = "  +3.69%  "
$ws.Range("E9").Value = "  +1.67%  "
$ws.Range("D10").Value = "'21.83"
$ws.Range("E10").Value = "  +5.67%  "
$ws.Range("D11").Value = "'0.07796"
$ws.Range("E11").Value = "  +1.84%  "
$ws.Range("D12").Value = "1.683.03"
$ws.Range("E12").Value = "  +2.85%  "
$ws.Range("D13").Value = "'4.513"
$ws.Range("E13").Value = "  +2.21%  "
$ws.Range("D14").Value = "'0.5585"
$ws.Range("E14").Value = "  +0.99%  "
$ws.Range("D15").Value = "0.0₅8323"
$ws.Range("E15").Value = "  +0.69%  "
$ws.Range("D16").Value = "'65.67"
$ws.Range("E16").Value = "  +1.12%  "
$ws.Range("D17").Value = "26.496.68"
$ws.Range("E17").Value = "  +1.79%  "
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "'4.783"
$ws.Range("E19").Value = "  +1.89%  "
$ws.Range("D20").Value = "'193.59"
$ws.Range("E20").Value = "  +2.99%  "
$ws.Range("E21").Value = "  +1.17%  "
$ws.Range("D22").Value = "'6.336"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'142.36"
$ws.Range("E24").Value = "  -2.01%  "
$ws.Range("D25").Value = "'0.1280"
$ws.Range("E25").Value = "  +5.15%  "
$ws.Range("D26").Value = "'7.415"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "'16.28"
$ws.Range("E27").Value = "  +3.08%  "
$ws.Range("D28").Value = "'1.446"
$ws.Range("E28").Value = "  +3.59%  "
$ws.Range("D29").Value = "'0.06279"
$ws.Range("E29").Value = "  +5.09%  "
$ws.Range("D30").Value = "'1.276"
$ws.Range("E30").Value = "  +1.76%  "
$ws.Range("D31").Value = "'3.614"
$ws.Range("E31").Value = "  +5.14%  "
$ws.Range("D32").Value = "'3.455"
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("D33").Value = "'1.689"
$ws.Range("E33").Value = "  +2.77%  "
$ws.Range("D34").Value = "'1.009"
$ws.Range("E34").Value = "  +2.46%  "
$ws.Range("D35").Value = "'0.6162"
$ws.Range("E35").Value = "  +8.41%  "
$ws.Range("E36").Value = "  +0.89%  "
$ws.Range("D38").Value = "'6.179"
$ws.Range("E38").Value = "  +7.25%  "
$ws.Range("D39").Value = "'0.01634"
$ws.Range("E39").Value = "  +1.18%  "
$ws.Range("D40").Value = "1.096.20"
$ws.Range("E40").Value = "  +5.95%  "
$ws.Range("E41").Value = "  +1.71%  "
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("D44").Value = "1.821.48"
$ws.Range("E44").Value = "  +1.98%  "
$ws.Range("D45").Value = "'57.89"
$ws.Range("E45").Value = "  +3.97%  "
$ws.Range("D46").Value = "'8.153"
$ws.Range("E46").Value = "  +1.16%  "
$ws.Range("D47").Value = "'0.9994"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("E48").Value = "  -2.94%  "
$ws.Range("D49").Value = "'0.05210"
$ws.Range("E49").Value = "  +1.00%  "
$ws.Range("E50").Value = "  +6.90%  "
$ws.Range("E51").Value = "  +2.17%  "
